# Generate Report for Handback
# - Mark both localized files as handed back (Status column) on the zh-cn and de-de sheets.
# - Record the "Latest Target File" and "Latest Handback File" (with hyperlinks) for each row.
# - Stamp the "Latest Handback DateTime" for each row.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status -> handed back
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (F) - mirrors the source markdown file name, with hyperlink
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7394dd430531af166018016b0f79eee11f0ad099/e2e/3b5cfbe7-c379-4594-aec4-2cf4c879c669.md", "", "", "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md")

# Latest Handback File (G) - the zh-cn xlf that was handed back, with hyperlink
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/133d6e38e32e8d84a534dea8574db6c3347de31b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.zh-cn.xlf", "", "", "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/7394dd430531af166018016b0f79eee11f0ad099/e2e/9e9b3e7d-d93e-447a-bbad-150e428577a6.md", "", "", "9e9b3e7d-d93e-447a-bbad-150e428577a6.md")

$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/133d6e38e32e8d84a534dea8574db6c3347de31b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.zh-cn.xlf", "", "", "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.zh-cn.xlf")

# Latest Handback DateTime (H) for both rows
$wsZh.Range("H2").Value = "2016-03-24 08:26:03"
$wsZh.Range("H3").Value = "2016-03-24 08:26:03"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status -> handed back
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Latest Target File (F) - mirrors the source markdown file name, with hyperlink
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7394dd430531af166018016b0f79eee11f0ad099/e2e/3b5cfbe7-c379-4594-aec4-2cf4c879c669.md", "", "", "3b5cfbe7-c379-4594-aec4-2cf4c879c669.md")

# Latest Handback File (G) - the de-de xlf that was handed back, with hyperlink
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d858dfe7d23838b9835e39ae9b10d2ec1a563bd1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.de-de.xlf", "", "", "3b5cfbe7-c379-4594-aec4-2cf4c879c669.474a2621abc4c5f035b11ee2fc489a3cb6e0da35.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/7394dd430531af166018016b0f79eee11f0ad099/e2e/9e9b3e7d-d93e-447a-bbad-150e428577a6.md", "", "", "9e9b3e7d-d93e-447a-bbad-150e428577a6.md")

$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d858dfe7d23838b9835e39ae9b10d2ec1a563bd1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.de-de.xlf", "", "", "9e9b3e7d-d93e-447a-bbad-150e428577a6.fcc536d0c36bd9bd6200ad580d82b1c5572e7520.de-de.xlf")

# Latest Handback DateTime (H) for both rows
$wsDe.Range("H2").Value = "2016-03-24 08:26:14"
$wsDe.Range("H3").Value = "2016-03-24 08:26:14"
